$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quiz/homework column header: G1 = "H06"
$ws.Range("G1").Value = "H06"

# New scores entered in column G (H06 grades) for rows 2-16
$ws.Range("G2").Value = 11
$ws.Range("G3").Value = 10
$ws.Range("G4").Value = 9.75
$ws.Range("G5").Value = 8.5
$ws.Range("G6").Value = 9.75
$ws.Range("G7").Value = 9
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 10.5
$ws.Range("G10").Value = 6.5
$ws.Range("G11").Value = 8
$ws.Range("G12").Value = 10.5
$ws.Range("G13").Value = 9
$ws.Range("G14").Value = 0
$ws.Range("G15").Value = 9.75
$ws.Range("G16").Value = 9

# Match the centered formatting already used by the rest of the row (G7:G10 region)
# for the newly filled-in G11:G17 cells
$ws.Range("G11:G17").HorizontalAlignment = -4108

# Corrected attendance scores for student 4 (row 5)
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 7

# Corrected attendance formula for student 9 (row 10)
$ws.Range("F10").Formula = "=(11/12)*10"

# Move the active selection to B23
$ws.Range("B23").Select()
